# New Local Strings added
#
# Adds two new localized string pairs (key + Turkish translation) to the
# tr.xlsx localization sheet, widens column A to fit the longer keys, and
# appends the two new rows with the same formatting as the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A so the longer resource keys are readable.
$ws.Columns.Item(1).ColumnWidth = 41.46

# Append the two new localization entries right after the last existing row.
$ws.Range("A96").Value = "YouAreNotAuthorizedToEditSharedExam"
$ws.Range("B96").Value = "Ortak deneme sinavlari sadece  TestOkur yoneticileri tarafindan guncellenebilir."
$ws.Range("A97").Value = "YouAreNotAuthorizedToDeleteSharedExam"
$ws.Range("B97").Value = "Ortak deneme sinavlari sadece  TestOkur yoneticileri tarafindan silinebilir."

# Match the formatting (style) already used by the rest of the data rows.
$ws.Range("A95:B95").Copy()
$ws.Range("A96:B97").PasteSpecial(-4122)
